$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.109.19"
$ws.Range("E2").Value = "  +1.95%  "

$ws.Range("D3").Value = "3.437.99"

$ws.Range("E4").Value = "  +0.15%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "409.49"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.98%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "129.52"
$c.ClearFormats()
$ws.Range("E6").Value = "  -2.65%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.628"
$c.ClearFormats()
$ws.Range("E7").Value = "  +6.63%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E8").Value = "  -0.06%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.756"
$c.ClearFormats()
$ws.Range("E9").Value = "  +13.07%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.141"
$c.ClearFormats()
$ws.Range("E10").Value = "  +16.24%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "43.33"
$c.ClearFormats()
$ws.Range("E11").Value = "  +2.66%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.140"
$c.ClearFormats()
$ws.Range("E12").Value = "  -0.37%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "8.80"
$c.ClearFormats()
$ws.Range("E13").Value = "  +6.02%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "20.32"
$c.ClearFormats()
$ws.Range("E14").Value = "  +3.96%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0000195"
$c.ClearFormats()
$ws.Range("E15").Value = "  +53.36%  "

$ws.Range("D16").Value = "3.375.84"
$ws.Range("E16").Value = "  +0.49%  "

$ws.Range("E17").Value = "  +3.18%  "

$ws.Range("D18").Value = "62.141.28"
$ws.Range("E18").Value = "  +1.98%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.43"
$c.ClearFormats()
$ws.Range("E19").Value = "  +3.43%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "373.14"
$c.ClearFormats()
$ws.Range("E20").Value = "  +22.06%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "87.57"
$c.ClearFormats()
$ws.Range("E21").Value = "  +4.80%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "3.18"
$c.ClearFormats()
$ws.Range("E22").Value = "  -1.01%  "

$ws.Range("E23").Value = "  +5.01%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "3.21"
$c.ClearFormats()
$ws.Range("E24").Value = "  +2.86%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "31.71"
$c.ClearFormats()
$ws.Range("E25").Value = "  +8.17%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "4.80"
$c.ClearFormats()
$ws.Range("E26").Value = "  +0.61%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.45"
$c.ClearFormats()
$ws.Range("E27").Value = "  +1.89%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.69"
$c.ClearFormats()
$ws.Range("E28").Value = "  +2.15%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.73"
$c.ClearFormats()
$ws.Range("E29").Value = "  +10.60%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "43.98"
$c.ClearFormats()
$ws.Range("E30").Value = "  +7.25%  "

$ws.Range("E31").Value = "  -0.48%  "

$ws.Range("E32").Value = "  -0.01%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "11.83"
$c.ClearFormats()
$ws.Range("E33").Value = "  +4.97%  "

$ws.Range("E34").Value = "  -0.03%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0494"
$c.ClearFormats()
$ws.Range("E35").Value = "  +3.13%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "52.20"
$c.ClearFormats()
$ws.Range("E36").Value = "  +0.95%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E37").Value = "  +0.18%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.36"
$c.ClearFormats()
$ws.Range("E38").Value = "  -1.73%  "

$ws.Range("E39").Value = "  +0.71%  "

$ws.Range("E40").Value = "  +6.94%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "143.35"
$c.ClearFormats()
$ws.Range("E41").Value = "  +4.52%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.312"
$c.ClearFormats()
$ws.Range("E42").Value = "  +9.29%  "

$ws.Range("E43").Value = "  -0.51%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "4.01"
$c.ClearFormats()
$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("E45").Value = "  +0.61%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "21.75"
$c.ClearFormats()
$ws.Range("E47").Value = "  +1.63%  "

$ws.Range("D48").Value = "2.112.36"
$ws.Range("E48").Value = "  -0.37%  "

$ws.Range("E49").Value = "  -0.22%  "

$ws.Range("E50").Value = "  +3.01%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0362"
$c.ClearFormats()
$ws.Range("E51").Value = "  +6.72%  "
